$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The three time-log rows (86-88) were reordered: the "Editing ERD" entry
# that used to be first (row 86) now sits last (row 88); the other two rows
# shift up by one.
#   old row 86 (Editing ERD / Kern Philip)          -> new row 88
#   old row 87 (Connect website & db / Kern Philip) -> new row 86
#   old row 88 (Creating final ERD / Sarvan Amel)   -> new row 87

# New row 86: "Trying to connect website and database" / Kern Philip
$ws.Range("A86").Value = $null
$ws.Range("B86").Value = "Trying to connect website and database"
$ws.Range("C86").Value = 1.5
$ws.Range("D86").Value = "Kern Philip"

# New row 87: "Creating final ERD" / Sarvan Amel (with the date that used
# to sit on the old row 88)
$ws.Range("A87").Value = 44591
$ws.Range("B87").Value = "Creating final ERD"
$ws.Range("C87").Value = 1.33
$ws.Range("D87").Value = "Sarvan Amel"

# New row 88: "Editing ERD" / Kern Philip
$ws.Range("A88").Value = $null
$ws.Range("B88").Value = "Editing ERD"
$ws.Range("C88").Value = 1
$ws.Range("D88").Value = "Kern Philip"

# Update the rolled-up totals so they keep referencing the right rows.
$ws.Range("B94").Formula = "=C76+C79+C86+C88+C89"
$ws.Range("B96").Formula = "=C91+C87+C83+C82+C81+C80+C77"

# Restore the active-cell selection recorded in the saved workbook.
$ws.Range("B95").Select()
